$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 41: "LL(k) Grammars (continued)" - change LL(1) to LL(k) in body
# ---------------------------------------------------------------------
$s41 = $p.Slides.Item(41)
$sh41 = $s41.Shapes.Item(4)
$tr41 = $sh41.TextFrame.TextRange
$para41 = $tr41.Paragraphs(5, 1)
$run41 = $para41.Runs(1, 1)
$run41.Text = "In practice, the syntax of most programming languages can be defined, or at least closely approximated, by an LL(k) grammar"

# ---------------------------------------------------------------------
# Slide 43: "Recursive Decent Parsing" - update sample code block
# ---------------------------------------------------------------------
$s43 = $p.Slides.Item(43)
$sh43 = $s43.Shapes.Item(2)

# Resize/reposition the content placeholder to match the edited text box
$sh43.Left = 36.12496062992126
$sh43.Width = 655.2

$tr43 = $sh43.TextFrame.TextRange

# --- Paragraph 2: "parseLoop()  // called when parsing the outer loop"
#     becomes a single run: "parseLoopStmt()  // called when parsing the outer loop"
$para2 = $tr43.Paragraphs(2, 1)
$deadChars2 = $para2.Characters(1, 9)     # "parseLoop"
$deadChars2.Delete()

$tr43 = $sh43.TextFrame.TextRange
$para2 = $tr43.Paragraphs(2, 1)
$run2 = $para2.Runs(1, 1)
$run2.Text = "parseLoopStmt()         // called when parsing the outer loop"

# --- Paragraph 6: "parseLoop()   // called when paring the inner loop"
#     becomes 3 runs:
#       "      parseLoopStmt()   // called " / "when parsing " / "the inner loop"
$tr43 = $sh43.TextFrame.TextRange
$para6 = $tr43.Paragraphs(6, 1)

# fix "paring" -> "parsing" first (chars 31-42 within the paragraph = "when paring ")
$midSel = $para6.Characters(31, 12)
$midSel.Text = "when parsing "

# rename parseLoop -> parseLoopStmt (run 2 of the paragraph)
$tr43 = $sh43.TextFrame.TextRange
$para6 = $tr43.Paragraphs(6, 1)
$run2of6 = $para6.Runs(2, 1)
$run2of6.Text = "parseLoopStmt"

# merge the leading "      ", "parseLoopStmt" runs into the 3rd run
$tr43 = $sh43.TextFrame.TextRange
$para6 = $tr43.Paragraphs(6, 1)
$leadChars = $para6.Characters(1, 19)     # "      parseLoopStmt"
$leadChars.Delete()

$tr43 = $sh43.TextFrame.TextRange
$para6 = $tr43.Paragraphs(6, 1)
$run1of6 = $para6.Runs(1, 1)
$run1of6.Text = "      parseLoopStmt()   // called "
